# The deck currently uses the "Integral" theme (ppt/theme/theme1.xml) for
# its slide master, while the default "Office Theme" colours
# (ppt/theme/theme2.xml) sit unused on the notes master.
#
# The authored edit swaps the two themes so the slide master picks up the
# stock "Office Theme" palette. We reproduce that by driving the
# presentation's 12-colour theme colour scheme (exposed on a Slide, and
# shared by every slide because they all follow the one slide master) to
# the standard Office theme RGB values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office theme colours (Dark1, Light1, Dark2, Light2, Accent1-6, Hyperlink,
# FollowedHyperlink), expressed as VBA RGB() long values (R + G*256 + B*65536).
$tcs.Colors(1).RGB  = 0        # Dark 1    -> 000000
$tcs.Colors(2).RGB  = 16777215 # Light 1   -> FFFFFF
$tcs.Colors(3).RGB  = 6968388  # Dark 2    -> 44546A
$tcs.Colors(4).RGB  = 15132391 # Light 2   -> E7E6E6
$tcs.Colors(5).RGB  = 13998939 # Accent 1  -> 5B9BD5
$tcs.Colors(6).RGB  = 3243501  # Accent 2  -> ED7D31
$tcs.Colors(7).RGB  = 10855845 # Accent 3  -> A5A5A5
$tcs.Colors(8).RGB  = 49407    # Accent 4  -> FFC000
$tcs.Colors(9).RGB  = 12874308 # Accent 5  -> 4472C4
$tcs.Colors(10).RGB = 4697456  # Accent 6  -> 70AD47
$tcs.Colors(11).RGB = 12673797 # Hyperlink -> 0563C1
$tcs.Colors(12).RGB = 7491477  # Followed Hyperlink -> 954F72
